# 自动更新Excel文件 - 2025-12-25 23:13:27
# Decrement the "剩余" (remaining days) value in column E by 1 for each data
# row (rows 2-99), except row 36 whose start date (column F) is malformed
# and therefore was not recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    $cell = $ws.Cells.Item($row, 5)   # Column E
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current - 1
    }
}
